$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1653225806451613
$ws.Range("C2").Value = 0.5887096774193549
$ws.Range("J2").Value = 0.0282258064516129
$ws.Range("P2").Value = 0.125
$ws.Range("S2").Value = 0.09274193548387097
$ws.Range("C3").Value = 0.03267973856209151
$ws.Range("J3").Value = 0.05228758169934641
$ws.Range("P3").Value = 0.7320261437908496
$ws.Range("S3").Value = 0.1830065359477124
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.7083333333333334
$ws.Range("S4").Value = 0.2291666666666667
$ws.Range("B6").Value = 0.04743083003952569
$ws.Range("D6").Value = 0.007905138339920948
$ws.Range("E6").Value = 0.003952569169960474
$ws.Range("F6").Value = 0.04347826086956522
$ws.Range("J6").Value = 0.1936758893280632
$ws.Range("O6").Value = 0.0158102766798419
$ws.Range("Q6").Value = 0.150197628458498
$ws.Range("R6").Value = 0.1067193675889328
$ws.Range("S6").Value = 0.4308300395256917
$ws.Range("B7").Value = 0.1067961165048544
$ws.Range("D7").Value = 0.01941747572815534
$ws.Range("E7").Value = 0.009708737864077669
$ws.Range("F7").Value = 0.07281553398058252
$ws.Range("J7").Value = 0.0970873786407767
$ws.Range("O7").Value = 0.01941747572815534
$ws.Range("Q7").Value = 0.1650485436893204
$ws.Range("R7").Value = 0.09223300970873786
$ws.Range("S7").Value = 0.4174757281553398
$ws.Range("B8").Value = 0.09147609147609148
$ws.Range("D8").Value = 0.01871101871101871
$ws.Range("E8").Value = 0.002079002079002079
$ws.Range("F8").Value = 0.07276507276507277
$ws.Range("J8").Value = 0.103950103950104
$ws.Range("O8").Value = 0.02494802494802495
$ws.Range("Q8").Value = 0.1683991683991684
$ws.Range("R8").Value = 0.1205821205821206
$ws.Range("S8").Value = 0.3970893970893971
$ws.Range("B9").Value = 0.08383233532934131
$ws.Range("D9").Value = 0.02395209580838323
$ws.Range("F9").Value = 0.0718562874251497
$ws.Range("J9").Value = 0.1077844311377246
$ws.Range("O9").Value = 0.01197604790419162
$ws.Range("Q9").Value = 0.1856287425149701
$ws.Range("R9").Value = 0.1017964071856287
$ws.Range("S9").Value = 0.4131736526946108
$ws.Range("B10").Value = 0.09083402146985962
$ws.Range("D10").Value = 0.02394715111478117
$ws.Range("F10").Value = 0.09083402146985962
$ws.Range("J10").Value = 0.1007431874483898
$ws.Range("O10").Value = 0.02394715111478117
$ws.Range("Q10").Value = 0.1816680429397192
$ws.Range("R10").Value = 0.1114781172584641
$ws.Range("S10").Value = 0.3765483071841453
$ws.Range("G11").Value = 0.1146496815286624
$ws.Range("J11").Value = 0.09554140127388536
$ws.Range("K11").Value = 0.1815286624203822
$ws.Range("L11").Value = 0.589171974522293
$ws.Range("S11").Value = 0.01910828025477707
$ws.Range("G12").Value = 0.7540106951871658
$ws.Range("J12").Value = 0.1978609625668449
$ws.Range("K12").Value = 0.0053475935828877
$ws.Range("L12").Value = 0.0106951871657754
$ws.Range("S12").Value = 0.03208556149732621
$ws.Range("G13").Value = 0.6938775510204082
$ws.Range("J13").Value = 0.2653061224489796
$ws.Range("S13").Value = 0.04081632653061224
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("F15").Value = 0.01171875
$ws.Range("H15").Value = 0.15234375
$ws.Range("J15").Value = 0.3359375
$ws.Range("K15").Value = 0.05859375
$ws.Range("M15").Value = 0.01953125
$ws.Range("N15").Value = 0.0078125
$ws.Range("O15").Value = 0.08984375
$ws.Range("S15").Value = 0.26171875
$ws.Range("F16").Value = 0.02285714285714286
$ws.Range("H16").Value = 0.1542857142857143
$ws.Range("I16").Value = 0.09142857142857143
$ws.Range("J16").Value = 0.36
$ws.Range("K16").Value = 0.1142857142857143
$ws.Range("M16").Value = 0.02285714285714286
$ws.Range("N16").Value = 0.01142857142857143
$ws.Range("O16").Value = 0.05142857142857143
$ws.Range("S16").Value = 0.1714285714285714
$ws.Range("F17").Value = 0.02743142144638404
$ws.Range("H17").Value = 0.2044887780548628
$ws.Range("I17").Value = 0.07231920199501247
$ws.Range("J17").Value = 0.4014962593516209
$ws.Range("K17").Value = 0.09725685785536159
$ws.Range("M17").Value = 0.009975062344139651
$ws.Range("O17").Value = 0.04488778054862843
$ws.Range("S17").Value = 0.14214463840399
$ws.Range("F18").Value = 0.02723735408560311
$ws.Range("H18").Value = 0.178988326848249
$ws.Range("I18").Value = 0.09727626459143969
$ws.Range("J18").Value = 0.4046692607003891
$ws.Range("K18").Value = 0.08171206225680934
$ws.Range("M18").Value = 0.02723735408560311
$ws.Range("O18").Value = 0.05836575875486381
$ws.Range("S18").Value = 0.1245136186770428
$ws.Range("F19").Value = 0.01386748844375963
$ws.Range("H19").Value = 0.224191063174114
$ws.Range("I19").Value = 0.06471494607087827
$ws.Range("J19").Value = 0.3528505392912172
$ws.Range("K19").Value = 0.1217257318952234
$ws.Range("M19").Value = 0.02311248073959938
$ws.Range("N19").Value = 0.0007704160246533128
$ws.Range("O19").Value = 0.08320493066255778
$ws.Range("S19").Value = 0.1155624036979969
